# Daily attendance processing - 2025-12-13 20:51:21
# Normalize the "Recorded By" (column G) lists so that any exact-case
# "System" entries are moved to the end of the comma-separated list,
# while preserving the relative order of all other entries (including a
# differently-cased "system").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Transform-RecordedBy($s) {
    $parts = $s.Split(",")
    $nonSystem = New-Object System.Collections.ArrayList
    $systemCount = 0
    foreach ($p in $parts) {
        $trimmedPart = $p.Trim()
        if ($trimmedPart.Equals("System")) {
            $systemCount++
        } else {
            $nonSystem.Add($trimmedPart) | Out-Null
        }
    }
    for ($i = 0; $i -lt $systemCount; $i++) {
        $nonSystem.Add("System") | Out-Null
    }
    return ($nonSystem -join ", ")
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($val -ne $null) {
        $newVal = Transform-RecordedBy $val
        if (-not $newVal.Equals($val)) {
            $cell.Value = $newVal
        }
    }
}
